$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '297.21'
Set-TextValue $ws.Range('E2') '1.63%'
Set-TextValue $ws.Range('E3') '3.51%'
Set-TextValue $ws.Range('D4') '5.004'
Set-TextValue $ws.Range('E4') '-0.21%'
Set-TextValue $ws.Range('D5') '0.07516'
Set-TextValue $ws.Range('E5') '2.69%'
$ws.Range('B6').Value = 'GateToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws.Range('D6') '4.376'
Set-TextValue $ws.Range('E6') '1.79%'
$ws.Range('B7').Value = 'FTXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws.Range('D7') '1.582'
Set-TextValue $ws.Range('E7') '3.49%'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D8') '0.9262'
Set-TextValue $ws.Range('E8') '-0.06%'
$ws.Range('B9').Value = 'BTSEToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws.Range('D9') '2.401'
Set-TextValue $ws.Range('E9') '1.23%'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws.Range('D10') '0.1197'
Set-TextValue $ws.Range('E10') '0.50%'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws.Range('D11') '0.1822'
Set-TextValue $ws.Range('E11') '4.69%'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws.Range('D12') '0.08946'
Set-TextValue $ws.Range('E12') '3.56%'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws.Range('D13') '0.04083'
Set-TextValue $ws.Range('E13') '-5.42%'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range('D14') '0.1048'
Set-TextValue $ws.Range('E14') '-0.43%'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range('D15') '0.001280'
Set-TextValue $ws.Range('E15') '0.01%'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws.Range('D16') '0.005811'
Set-TextValue $ws.Range('E16') '-3.20%'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range('D17') '3.357'
Set-TextValue $ws.Range('E17') '0.46%'
Set-TextValue $ws.Range('D18') '0.3314'
Set-TextValue $ws.Range('E18') '0.75%'
Set-TextValue $ws.Range('D19') '8.118'
Set-TextValue $ws.Range('E19') '1.79%'
Set-TextValue $ws.Range('D20') '0.1392'
Set-TextValue $ws.Range('E20') '0.10%'
Set-TextValue $ws.Range('E21') '11.01%'
Set-TextValue $ws.Range('D22') '0.04096'
Set-TextValue $ws.Range('E22') '3.90%'
Set-TextValue $ws.Range('E23') '0.44%'
Set-TextValue $ws.Range('D24') '0.003897'
Set-TextValue $ws.Range('E24') '3.15%'
Set-TextValue $ws.Range('E25') '-3.92%'
Set-TextValue $ws.Range('D38') '0.02404'
Set-TextValue $ws.Range('E38') '5.27%'
Set-TextValue $ws.Range('D39') '0.05205'
Set-TextValue $ws.Range('E39') '4.67%'
Set-TextValue $ws.Range('D40') '0.006306'
Set-TextValue $ws.Range('E40') '11.95%'
Set-TextValue $ws.Range('D41') '0.007834'
Set-TextValue $ws.Range('E41') '1.98%'
Set-TextValue $ws.Range('E42') '3.24%'
Set-TextValue $ws.Range('D43') '0.007413'
Set-TextValue $ws.Range('E43') '1.04%'
Set-TextValue $ws.Range('D44') '0.007270'
Set-TextValue $ws.Range('E44') '-0.07%'
Set-TextValue $ws.Range('D45') '0.2964'
Set-TextValue $ws.Range('E45') '1.35%'
Set-TextValue $ws.Range('D46') '0.00006589'
Set-TextValue $ws.Range('E46') '4.20%'
Set-TextValue $ws.Range('E47') '0.01%'
Set-TextValue $ws.Range('D48') '0.03158'
Set-TextValue $ws.Range('E48') '48.07%'
Set-TextValue $ws.Range('D49') '0.004204'
Set-TextValue $ws.Range('E49') '0.05%'
Set-TextValue $ws.Range('E50') '0.01%'
Set-TextValue $ws.Range('E51') '0.01%'
